# Auto-generated update script for horarios-141-2026-01-22.xlsx
# Applies the scraper refresh: new 'Ultima actualizacion' / 'Total filas' headers
# plus updated / newly-appended schedule rows across all three sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range("A2").Value = 'Última actualización: 06:24:49'
$ws.Range("A3").Value = 'Total filas: 63'

$rows = @(
    @(37, '06:24:49', '06:29', '86_EST CHICA-ESC AGRARIA', 5, 'LP1912'),
    @(39, '06:24:49', '06:31', '16_SANTA ANA', 7, 'LP1912'),
    @(40, '06:24:49', '06:32', '23_HERNANDEZ', 8, 'LP1912'),
    @(41, '06:24:49', '06:44', '225_C ROCA-H SUR', 20, 'LP1912'),
    @(42, '06:24:49', '06:46', '215C_EL PATO', 22, 'LP1912'),
    @(43, '05:23:05', '06:47', '215C_EL PATO', 84, 'LP1912'),
    @(44, '06:24:49', '06:59', '14_ABASTO', 35, 'LP1912'),
    @(45, '05:23:05', '07:00', '14_ABASTO', 97, 'LP1912'),
    @(46, '06:24:49', '07:01', '16_SANTA ANA', 37, 'LP1912'),
    @(47, '05:54:50', '07:04', '23_HERNANDEZ', 70, 'LP1912'),
    @(48, '06:24:49', '07:05', '15_ABASTO', 41, 'LP1912'),
    @(49, '06:24:49', '07:06', '225_GOMEZ', 42, 'LP1912'),
    @(50, '05:23:05', '07:07', '225_GOMEZ', 104, 'LP1912'),
    @(51, '06:24:49', '07:11', '215A_EL PATO', 47, 'LP1912'),
    @(52, '05:23:05', '07:12', '215A_EL PATO', 109, 'LP1912'),
    @(53, '06:24:49', '07:15', '11_ETCHEVERRY', 51, 'LP1912'),
    @(54, '05:23:05', '07:16', '11_ETCHEVERRY', 113, 'LP1912'),
    @(55, '06:24:49', '07:21', '26_HERNANDEZ', 57, 'LP1912'),
    @(56, '06:24:49', '07:23', '10_OLMOS', 59, 'LP1912'),
    @(57, '06:24:49', '07:31', '11_ETCHEVERRY', 67, 'LP1912'),
    @(58, '06:24:49', '07:32', '84_COLONIA URQUIZA-ESC 49', 68, 'LP1912'),
    @(59, '06:24:49', '07:36', '27_EL RETIRO', 72, 'LP1912'),
    @(60, '06:24:49', '07:39', '10_OLMOS', 75, 'LP1912'),
    @(61, '05:54:50', '07:46', '16_SANTA ANA', 112, 'LP1912'),
    @(62, '06:24:49', '07:47', '14_ABASTO', 83, 'LP1912'),
    @(63, '06:24:49', '07:51', '215D_EL PATO', 87, 'LP1912'),
    @(64, '06:24:49', '08:05', '23_HERNANDEZ', 101, 'LP1912'),
    @(65, '06:24:49', '08:12', '15_ABASTO', 108, 'LP1912'),
    @(66, '06:24:49', '08:21', '26_HERNANDEZ', 117, 'LP1912'),
    @(67, '06:24:49', '08:22', '16_P MOR-SANTA ANA', 118, 'LP1912'),
    @(68, '06:24:49', '08:23', '215B_EL PATO', 119, 'LP1912')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range("A2").Value = 'Última actualización: 06:24:49'
$ws.Range("A3").Value = 'Total filas: 17'

$rows = @(
    @(17, '06:24:49', '06:46', '215C_EL PATO', 22, 'LP1912'),
    @(19, '06:24:49', '07:11', '215A_EL PATO', 47, 'LP1912'),
    @(21, '06:24:49', '07:51', '215D_EL PATO', 87, 'LP1912'),
    @(22, '06:24:49', '08:23', '215B_EL PATO', 119, 'LP1912')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range("A2").Value = 'Última actualización: 06:24:49'
$ws.Range("A3").Value = 'Total filas: 10'

$rows = @(
    @(10, '06:24:49', '06:32', '215C_LA PLATA', 8, 'L6203'),
    @(12, '06:24:49', '06:59', '215B_LP-P MOR-1 Y 57', 35, 'L6173'),
    @(14, '06:24:49', '07:35', '215A_LA PLATA', 71, 'L6173'),
    @(15, '06:24:49', '08:06', '215C_LA PLATA', 102, 'L6203')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
